# Update the "სოციალური პაკეტის მიმღებთა რიცხოვნობა" figures for
# ოზურგეთი municipality (row 4, years 2015-2021 / columns E-K) with the
# corrected/refreshed values from the source ministry data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 2703   # 2015
$ws.Range("F4").Value = 2722   # 2016
$ws.Range("G4").Value = 2702   # 2017
$ws.Range("H4").Value = 2697   # 2018
$ws.Range("I4").Value = 2830   # 2019
$ws.Range("J4").Value = 2854   # 2020
$ws.Range("K4").Value = 2879   # 2021
